$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.112.20"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.874.01"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'312.89"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").Value = "'0.9990"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "'0.5053"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.3843"
$ws.Range("E8").Value = "  -2.03%  "
$ws.Range("D9").Value = "'0.09029"
$ws.Range("E9").Value = "  -5.69%  "
$ws.Range("D10").Value = "'1.121"
$ws.Range("E10").Value = "  -1.20%  "
$ws.Range("D11").Value = "'41.70"
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").Value = "'6.374"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "'20.80"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "1.874.44"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").Value = "'7.266"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "'0.00001108"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").Value = "'91.33"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "'0.06652"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("E20").Value = "  +2.29%  "
$ws.Range("D21").Value = "'0.9992"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").Value = "'6.138"
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("D23").Value = "28.111.71"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'11.49"
$ws.Range("E24").Value = "  +2.02%  "
$ws.Range("D25").Value = "'2.251"
$ws.Range("E25").Value = "  -2.53%  "
$ws.Range("D26").Value = "2.085.65"
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("D27").Value = "'2.538"
$ws.Range("E27").Value = "  -4.61%  "
$ws.Range("D28").Value = "'20.89"
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("D29").Value = "'157.16"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'127.79"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("D31").Value = "'0.1065"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").Value = "'1.065"
$ws.Range("E32").Value = "  -1.48%  "
$ws.Range("D33").Value = "'5.612"
$ws.Range("D34").Value = "'3.592"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").Value = "'9.441"
$ws.Range("E35").Value = "  -1.50%  "
$ws.Range("D36").Value = "'0.06587"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "'0.02405"
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("D38").Value = "'0.2195"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("D39").Value = "'1.291"
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("D40").Value = "'1.212"
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("D41").Value = "'0.6400"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("D43").Value = "'4.928"
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("D44").Value = "'0.9991"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'13.27"
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.6029"
$ws.Range("E46").Value = "  +0.89%  "
$ws.Range("D47").Value = "'1.276"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("D49").Value = "'1.238"
$ws.Range("E49").Value = "  +4.68%  "
$ws.Range("D50").Value = "'2.002"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").Value = "'121.20"
$ws.Range("E51").Value = "  -1.35%  "
